# Fill in the actual worked hours for the week of 03/07/17 on the "7 Days"
# sheet.  Only the raw time-in / time-out punches (and the "Expect. Hrs"
# column) are edited here -- every other cell on this sheet (H/L/M) is a
# formula and recalculates automatically, as do the rollup cells on the
# totals sheet ("Sheet1", e.g. B8 = '7 Days'!H66).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("7 Days")

# Monday (row 59): logged out at 17:00 instead of 16:00, expected 8 hrs
$ws.Range("E59").Value = 0.70833333333333337
$ws.Range("K59").Value = 8

# Tuesday (row 60): 09:15 - 12:00, 12:30 - 17:45
$ws.Range("B60").Value = 0.38541666666666669
$ws.Range("C60").Value = 0.5
$ws.Range("D60").Value = 0.52083333333333337
$ws.Range("E60").Value = 0.73958333333333337

# Wednesday (row 61): 09:00 - 12:00, 12:30 - 18:00
$ws.Range("B61").Value = 0.375
$ws.Range("C61").Value = 0.5
$ws.Range("D61").Value = 0.52083333333333337
$ws.Range("E61").Value = 0.75

# Thursday (row 62): 09:00 - 12:00, 12:30 - 18:00
$ws.Range("B62").Value = 0.375
$ws.Range("C62").Value = 0.5
$ws.Range("D62").Value = 0.52083333333333337
$ws.Range("E62").Value = 0.75

# Friday (row 63): 09:00 - 12:00, 12:30 - 18:00
$ws.Range("B63").Value = 0.375
$ws.Range("C63").Value = 0.5
$ws.Range("D63").Value = 0.52083333333333337
$ws.Range("E63").Value = 0.75

# Saturday (row 64): not expected to work
$ws.Range("K64").Value = 0

$excel.CalculateFullRebuild()
